$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new C-column value (logistic_tfidf), or $null if unchanged
$rowUpdates = @{
    7 = $null
    14 = 0.501
    21 = 0.344
    28 = 0.476
    35 = 0.388
    42 = 0.388
    49 = 0.476
    56 = 0.484
    63 = 0.36
    70 = 0.377
    77 = 0.476
    84 = 0.52
    91 = 0.335
    98 = 0.497
    105 = 0.442
    112 = 0.421
    119 = 0.422
    126 = 0.003
    133 = 0.001
    140 = $null
    147 = 0.002
    154 = 0.004
    161 = 0.004
    168 = 0.002
    175 = 0.003
    182 = 0.006
    189 = 0.003
    196 = 0.002
    203 = 0.003
    210 = 0.004
    217 = 0.002
    224 = 0.007
    231 = $null
    238 = $null
}

foreach ($row in $rowUpdates.Keys) {
    # Column B: n_sample interval label, change text "3000" -> "5000"
    $ws.Range("B$row").Value = "'5000"
    $ws.Range("B$row").Style = "Normal"

    # Column C: logistic_tfidf metric value, only update where diff specifies a new value
    $newC = $rowUpdates[$row]
    if ($null -ne $newC) {
        $ws.Range("C$row").Value = $newC
    }
}

Write-Host "Applied cap-sotu 3000->5000 interval fix to $($rowUpdates.Count) rows"
